# Updated symbol list on Mon Feb 13 18:33:34 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scraped quote snapshot. Values are kept as literal text (matching
# the sheet's existing inline-string convention for these columns) rather
# than being auto-coerced to numbers/percentages by Excel's type-sniffing,
# so each cell's NumberFormat is forced to Text ("@") immediately before
# the new literal is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "288.34" }
    @{ Cell = "E2"; Value = "-9.85%" }
    @{ Cell = "D3"; Value = "40.25" }
    @{ Cell = "E3"; Value = "-2.73%" }
    @{ Cell = "D4"; Value = "5.037" }
    @{ Cell = "E4"; Value = "-3.97%" }
    @{ Cell = "D5"; Value = "0.07282" }
    @{ Cell = "E5"; Value = "-5.81%" }
    @{ Cell = "D6"; Value = "4.281" }
    @{ Cell = "E6"; Value = "-1.30%" }
    @{ Cell = "D7"; Value = "1.514" }
    @{ Cell = "E7"; Value = "-10.88%" }
    @{ Cell = "D8"; Value = "0.9142" }
    @{ Cell = "E8"; Value = "-3.13%" }
    @{ Cell = "D9"; Value = "0.1193" }
    @{ Cell = "E9"; Value = "-3.66%" }
    @{ Cell = "D10"; Value = "0.1710" }
    @{ Cell = "E10"; Value = "-6.56%" }
    @{ Cell = "D11"; Value = "0.08620" }
    @{ Cell = "E11"; Value = "-6.08%" }
    @{ Cell = "D12"; Value = "0.04178" }
    @{ Cell = "E12"; Value = "-3.70%" }
    @{ Cell = "E13"; Value = "0.28%" }
    @{ Cell = "D14"; Value = "0.001274" }
    @{ Cell = "E14"; Value = "-1.50%" }
    @{ Cell = "D15"; Value = "0.005936" }
    @{ Cell = "E15"; Value = "-1.44%" }
    @{ Cell = "D16"; Value = "3.402" }
    @{ Cell = "E16"; Value = "1.84%" }
    @{ Cell = "D18"; Value = "0.3261" }
    @{ Cell = "E18"; Value = "-2.90%" }
    @{ Cell = "D19"; Value = "7.790" }
    @{ Cell = "E19"; Value = "1.01%" }
    @{ Cell = "D20"; Value = "0.1352" }
    @{ Cell = "E20"; Value = "-0.10%" }
    @{ Cell = "D21"; Value = "0.2885" }
    @{ Cell = "E21"; Value = "2.08%" }
    @{ Cell = "D22"; Value = "0.03854" }
    @{ Cell = "E22"; Value = "-4.39%" }
    @{ Cell = "D23"; Value = "0.001268" }
    @{ Cell = "E23"; Value = "0.19%" }
    @{ Cell = "D24"; Value = "0.003806" }
    @{ Cell = "E24"; Value = "-7.43%" }
    @{ Cell = "D25"; Value = "0.0001281" }
    @{ Cell = "E25"; Value = "0.70%" }
    @{ Cell = "D26"; Value = "0.0003726" }
    @{ Cell = "D38"; Value = "0.02298" }
    @{ Cell = "E38"; Value = "-9.70%" }
    @{ Cell = "D39"; Value = "0.04948" }
    @{ Cell = "E39"; Value = "-7.40%" }
    @{ Cell = "D40"; Value = "0.007102" }
    @{ Cell = "E40"; Value = "256.66%" }
    @{ Cell = "D41"; Value = "0.007698" }
    @{ Cell = "E41"; Value = "-1.01%" }
    @{ Cell = "D42"; Value = "0.1269" }
    @{ Cell = "E42"; Value = "-3.65%" }
    @{ Cell = "D43"; Value = "0.007365" }
    @{ Cell = "E43"; Value = "-0.01%" }
    @{ Cell = "D44"; Value = "0.007604" }
    @{ Cell = "E44"; Value = "-9.21%" }
    @{ Cell = "D45"; Value = "0.3118" }
    @{ Cell = "E45"; Value = "-1.90%" }
    @{ Cell = "D46"; Value = "0.00006364" }
    @{ Cell = "E46"; Value = "-5.13%" }
    @{ Cell = "D47"; Value = "0.00000000750" }
    @{ Cell = "E47"; Value = "-0.08%" }
    @{ Cell = "D48"; Value = "0.2414" }
    @{ Cell = "E48"; Value = "22.65%" }
    @{ Cell = "E49"; Value = "-0.11%" }
    @{ Cell = "D50"; Value = "0.00002101" }
    @{ Cell = "E50"; Value = "-0.08%" }
    @{ Cell = "D51"; Value = "0.0002001" }
    @{ Cell = "E51"; Value = "-0.08%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

